$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "主力流入1408亿"
$ws.Range("A3").Value = "主力流出1514亿"
$ws.Range("A4").Value = "主力净流入-106.2亿"
$ws.Range("A5").Value = "超大单462.8亿525.6亿"
$ws.Range("A6").Value = "大单945.2亿988.5亿"
$ws.Range("A7").Value = "中单1350亿1344亿"
$ws.Range("A8").Value = "小单1293亿1193亿"
